$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '36.418.44'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -2.76%  '
$c.Style = "Normal"

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.979.12'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -3.25%  '
$c.Style = "Normal"

# Row 4
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '245.38'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -2.70%  '
$c.Style = "Normal"

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.623'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -4.17%  '
$c.Style = "Normal"

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '58.77'
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -10.01%  '
$c.Style = "Normal"

# Row 8
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.374'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -7.52%  '
$c.Style = "Normal"

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '56.66'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -5.06%  '
$c.Style = "Normal"

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0875'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +10.13%  '
$c.Style = "Normal"

# Row 12
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '
$c.Style = "Normal"

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.854'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -6.92%  '
$c.Style = "Normal"

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '22.08'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -6.26%  '
$c.Style = "Normal"

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.271.71'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -3.14%  '
$c.Style = "Normal"

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '13.71'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -7.56%  '
$c.Style = "Normal"

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.45'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -4.93%  '
$c.Style = "Normal"

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.964.53'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -3.94%  '
$c.Style = "Normal"

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '36.302.88'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -2.67%  '
$c.Style = "Normal"

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0906'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +2.55%  '
$c.Style = "Normal"

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '70.36'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -4.09%  '
$c.Style = "Normal"

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.26'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -4.24%  '
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '234.58'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -2.06%  '
$c.Style = "Normal"

# Row 24
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.49'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -3.93%  '
$c.Style = "Normal"

# Row 26
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -2.81%  '
$c.Style = "Normal"

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.78'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -2.03%  '
$c.Style = "Normal"

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '164.88'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +1.93%  '
$c.Style = "Normal"

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.92'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -0.46%  '
$c.Style = "Normal"

# Row 30
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -1.53%  '
$c.Style = "Normal"

# Row 31
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -2.33%  '
$c.Style = "Normal"

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -0.45%  '
$c.Style = "Normal"

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.87'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -5.48%  '
$c.Style = "Normal"

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0647'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +3.13%  '
$c.Style = "Normal"

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.41'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -5.56%  '
$c.Style = "Normal"

# Row 36
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -4.67%  '
$c.Style = "Normal"

# Row 38
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -1.89%  '
$c.Style = "Normal"

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -7.04%  '
$c.Style = "Normal"

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.92'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -2.29%  '
$c.Style = "Normal"

# Row 41
$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c.Style = "Normal"
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.0964'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -5.19%  '
$c.Style = "Normal"

# Row 42
$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = 'TrustWalletToken'
$c.Style = "Normal"
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.21'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -6.59%  '
$c.Style = "Normal"

# Row 43
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -5.21%  '
$c.Style = "Normal"

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0213'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -2.62%  '
$c.Style = "Normal"

# Row 45
$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c.Style = "Normal"
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.08'
$c.Style = "Normal"

# Row 46
$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = 'InjectiveProtocol'
$c.Style = "Normal"
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '16.16'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -7.29%  '
$c.Style = "Normal"

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '91.04'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -4.79%  '
$c.Style = "Normal"

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.361.98'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -2.77%  '
$c.Style = "Normal"

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.42'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -5.30%  '
$c.Style = "Normal"

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.83'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -2.74%  '
$c.Style = "Normal"

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '45.28'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -4.14%  '
$c.Style = "Normal"

